$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 216.5
$ws.Range("I4").Value = 159.8
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 159.8
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -45.80000000000001
$ws.Range("N4").Value = -728

$ws.Range("H62").Value = 5476.727
$ws.Range("I62").Value = 4035.7144
$ws.Range("J62").Value = 7998.5
$ws.Range("K62").Value = 4035.7144
$ws.Range("L62").Value = 7998.5
$ws.Range("M62").Value = -3411.7144
$ws.Range("N62").Value = -9246.5

$ws.Range("H65").Value = 5476.727
$ws.Range("I65").Value = 4035.7144
$ws.Range("J65").Value = 7998.5
$ws.Range("K65").Value = 20178.572
$ws.Range("L65").Value = 39992.5
$ws.Range("M65").Value = -17058.572
$ws.Range("N65").Value = -46232.5

$ws.Range("H112").Value = 2310.7144
$ws.Range("J112").Value = 2516.6667
$ws.Range("L112").Value = 7550.000100000001
$ws.Range("N112").Value = -9766.000100000001

$ws.Range("H135").Value = 679.2
$ws.Range("I135").Value = 679.2
$ws.Range("K135").Value = 6112.8
$ws.Range("M135").Value = -3577.8

$ws.Range("H137").Value = 1459
$ws.Range("I137").Value = 1153
$ws.Range("K137").Value = 3459
$ws.Range("M137").Value = -909

$ws.Range("H138").Value = 1847.375
$ws.Range("J138").Value = 2750
$ws.Range("L138").Value = 8250
$ws.Range("N138").Value = -18530

$ws.Range("H141").Value = 2069.5833
$ws.Range("I141").Value = 1833.5
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 5500.5
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = -320.5
$ws.Range("N141").Value = -20110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 507.9091
$ws.Range("I97").Value = 486.44446
$ws.Range("J97").Value = 604.5
$ws.Range("K97").Value = 486.44446
$ws.Range("L97").Value = 604.5
$ws.Range("M97").Value = 9.555540000000008
$ws.Range("N97").Value = -1596.5

$ws.Range("H132").Value = 2939.5625
$ws.Range("I132").Value = 2939.5625
$ws.Range("K132").Value = 8818.6875
$ws.Range("M132").Value = -6288.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1651
$ws.Range("I5").Value = 1199.6666
$ws.Range("J5").Value = 3005
$ws.Range("K5").Value = 1199.6666
$ws.Range("L5").Value = 3005
$ws.Range("M5").Value = -1086.6666
$ws.Range("N5").Value = -3231

$ws.Range("H20").Value = 965.2727
$ws.Range("I20").Value = 795.5
$ws.Range("J20").Value = 1169
$ws.Range("K20").Value = 795.5
$ws.Range("L20").Value = 1169
$ws.Range("M20").Value = -548.5
$ws.Range("N20").Value = -1663

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H64").Value = 285.75
$ws.Range("J64").Value = 224
$ws.Range("L64").Value = 224
$ws.Range("N64").Value = -674

$ws.Range("H67").Value = 285.75
$ws.Range("J67").Value = 224
$ws.Range("L67").Value = 224
$ws.Range("N67").Value = -1784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 45037
$ws.Range("J45").Value = 45037
$ws.Range("L45").Value = 45037
$ws.Range("N45").Value = -46223

$ws.Range("H107").Value = 1456.7142
$ws.Range("I107").Value = 1550
$ws.Range("J107").Value = 1332.3334
$ws.Range("K107").Value = 1550
$ws.Range("L107").Value = 1332.3334
$ws.Range("M107").Value = 370
$ws.Range("N107").Value = -5172.3334

$ws.Range("H132").Value = 4199.8
$ws.Range("I132").Value = 4199.8
$ws.Range("K132").Value = 12599.4
$ws.Range("M132").Value = -10069.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 56.736843
$ws.Range("I2").Value = 70.333336
$ws.Range("K2").Value = 422.000016
$ws.Range("M2").Value = -309.000016

$ws.Range("H5").Value = 754.375
$ws.Range("I5").Value = 827.5
$ws.Range("K5").Value = 2482.5
$ws.Range("M5").Value = -2370.5

$ws.Range("H129").Value = 679.5
$ws.Range("J129").Value = 329
$ws.Range("L129").Value = 987
$ws.Range("N129").Value = -10987

$ws.Range("H131").Value = 1949.75
$ws.Range("I131").Value = 900
$ws.Range("K131").Value = 2700
$ws.Range("M131").Value = 2340

$ws.Range("H135").Value = 754.375
$ws.Range("I135").Value = 827.5
$ws.Range("K135").Value = 7447.5
$ws.Range("M135").Value = -4912.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3225.5
$ws.Range("I41").Value = 751
$ws.Range("J41").Value = 5700
$ws.Range("K41").Value = 751
$ws.Range("L41").Value = 5700
$ws.Range("M41").Value = -396
$ws.Range("N41").Value = -6410

$ws.Range("H70").Value = 166667660
$ws.Range("I70").Value = 166667660
$ws.Range("K70").Value = 166667660
$ws.Range("M70").Value = -166667390

$ws.Range("H73").Value = 166667660
$ws.Range("I73").Value = 166667660
$ws.Range("K73").Value = 166667660
$ws.Range("M73").Value = -166666724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 495.14285
$ws.Range("I16").Value = 411
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 411
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -241
$ws.Range("N16").Value = -1340

$ws.Range("H55").Value = 205.33333
$ws.Range("J55").Value = 220
$ws.Range("L55").Value = 220
$ws.Range("N55").Value = -566

$ws.Range("H68").Value = 1858
$ws.Range("J68").Value = 2596.3333
$ws.Range("L68").Value = 2596.3333
$ws.Range("N68").Value = -4094.3333

$ws.Range("H71").Value = 1858
$ws.Range("J71").Value = 2596.3333
$ws.Range("L71").Value = 12981.6665
$ws.Range("N71").Value = -20469.6665

$ws.Range("H82").Value = 1114.6666
$ws.Range("I82").Value = 1131.3334
$ws.Range("K82").Value = 1131.3334
$ws.Range("M82").Value = -770.3334

$ws.Range("H85").Value = 1114.6666
$ws.Range("I85").Value = 1131.3334
$ws.Range("K85").Value = 1131.3334
$ws.Range("M85").Value = 116.6666

$ws.Range("H93").Value = 1279.4
$ws.Range("I93").Value = 1249.25
$ws.Range("K93").Value = 1249.25
$ws.Range("M93").Value = -1.25

$ws.Range("H132").Value = 5749.5
$ws.Range("I132").Value = 4332.6665
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 12997.9995
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -10467.9995
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 47499.5
$ws.Range("I26").Value = 20000
$ws.Range("J26").Value = 74999
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 74999
$ws.Range("M26").Value = -19707
$ws.Range("N26").Value = -75585

$ws.Range("H132").Value = 2474.3
$ws.Range("I132").Value = 2320.7144
$ws.Range("J132").Value = 2832.6667
$ws.Range("K132").Value = 6962.1432
$ws.Range("L132").Value = 8498.000100000001
$ws.Range("M132").Value = -4432.1432
$ws.Range("N132").Value = -13558.0001
